# The workbook stores one row per fixture. Several fixture pairs had their
# odds/result columns (B, E:AD) swapped between the two physical rows - the
# id in column A (and the shared league/date in C:D) stay with the row, but
# every other value needs to move to the other row of the pair.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$swaps = @(
    @{
        RowA = 19
        RowB = 20
        ValuesForRowA = @{ 'B' = 7032917; 'E' = "FK Backa Topola"; 'F' = "FK Radnicki 1923"; 'G' = 1; 'H' = 0; 'I' = 0; 'J' = 0; 'K' = "H"; 'L' = 1.5; 'M' = 3.75; 'N' = 6.5; 'O' = 1.444; 'P' = 4; 'Q' = 6.5; 'R' = -1.25; 'S' = 1.975; 'T' = 1.825; 'U' = 2.75; 'V' = 1.95; 'W' = 1.85; 'X' = 0.444; 'Y' = -1; 'Z' = -1; 'AA' = -0.5; 'AB' = 0.4125; 'AC' = -1; 'AD' = 0.8500000000000001 }
        ValuesForRowB = @{ 'B' = 7032914; 'E' = "FK Vozdovac"; 'F' = "FK Radnik Surdulica"; 'G' = 1; 'H' = 1; 'I' = 0; 'J' = 1; 'K' = "D"; 'L' = 2.2; 'M' = 3.1; 'N' = 3.2; 'O' = 2.05; 'P' = 3.1; 'Q' = 3.5; 'R' = -0.25; 'S' = 1.75; 'T' = 2.05; 'U' = 2; 'V' = 1.775; 'W' = 2.025; 'X' = -1; 'Y' = 2.1; 'Z' = -1; 'AA' = -0.5; 'AB' = 0.5249999999999999; 'AC' = 0; 'AD' = 0 }
    },
    @{
        RowA = 31
        RowB = 32
        ValuesForRowA = @{ 'B' = 6979427; 'E' = "Javor Ivanjica"; 'F' = "FK Backa Topola"; 'G' = 0; 'H' = 3; 'I' = 0; 'J' = 3; 'K' = "A"; 'L' = 5.25; 'M' = 4; 'N' = 1.533; 'O' = 4; 'P' = 3.3; 'Q' = 1.85; 'R' = 0.5; 'S' = 1.875; 'T' = 1.925; 'U' = 2.25; 'V' = 1.825; 'W' = 1.975; 'X' = -1; 'Y' = -1; 'Z' = 0.8500000000000001; 'AA' = -1; 'AB' = 0.925; 'AC' = 0.825; 'AD' = -1 }
        ValuesForRowB = @{ 'B' = 6979431; 'E' = "FK Napredak"; 'F' = "Mladost Lucani"; 'G' = 0; 'H' = 1; 'I' = 0; 'J' = 0; 'K' = "A"; 'L' = 1.7; 'M' = 3.75; 'N' = 4.5; 'O' = 1.909; 'P' = 3.1; 'Q' = 4; 'R' = -0.5; 'S' = 1.975; 'T' = 1.825; 'U' = 2; 'V' = 1.875; 'W' = 1.925; 'X' = -1; 'Y' = -1; 'Z' = 3; 'AA' = -1; 'AB' = 0.825; 'AC' = -1; 'AD' = 0.925 }
    },
    @{
        RowA = 38
        RowB = 39
        ValuesForRowA = @{ 'B' = 6979435; 'E' = "FK Cukaricki"; 'F' = "FK Radnik Surdulica"; 'G' = 0; 'H' = 0; 'I' = 0; 'J' = 0; 'K' = "D"; 'L' = 1.4; 'M' = 3.75; 'N' = 7.5; 'O' = 1.5; 'P' = 3.6; 'Q' = 6; 'R' = -1; 'S' = 1.85; 'T' = 1.95; 'U' = 2.5; 'V' = 2; 'W' = 1.8; 'X' = -1; 'Y' = 2.6; 'Z' = -1; 'AA' = -1; 'AB' = 0.95; 'AC' = -1; 'AD' = 0.8 }
        ValuesForRowB = @{ 'B' = 6979433; 'E' = "Spartak Subotica"; 'F' = "FK Napredak"; 'G' = 1; 'H' = 3; 'I' = 0; 'J' = 2; 'K' = "A"; 'L' = 2.25; 'M' = 3.1; 'N' = 2.875; 'O' = 1.95; 'P' = 3; 'Q' = 3.75; 'R' = -0.5; 'S' = 2; 'T' = 1.8; 'U' = 2; 'V' = 1.85; 'W' = 1.95; 'X' = -1; 'Y' = -1; 'Z' = 2.75; 'AA' = -1; 'AB' = 0.8; 'AC' = 0.8500000000000001; 'AD' = -1 }
    },
    @{
        RowA = 90
        RowB = 91
        ValuesForRowA = @{ 'B' = 6979491; 'E' = "Radnicki Nis"; 'F' = "Spartak Subotica"; 'G' = 1; 'H' = 1; 'I' = 0; 'J' = 1; 'K' = "D"; 'L' = 1.95; 'M' = 3.25; 'N' = 3.7; 'O' = 1.65; 'P' = 3.5; 'Q' = 5; 'R' = -0.75; 'S' = 1.825; 'T' = 1.975; 'U' = 2.5; 'V' = 2; 'W' = 1.8; 'X' = -1; 'Y' = 2.5; 'Z' = -1; 'AA' = -1; 'AB' = 0.9750000000000001; 'AC' = -1; 'AD' = 0.8 }
        ValuesForRowB = @{ 'B' = 6978747; 'E' = "IMT Novi Belgrade"; 'F' = "Red Star Belgrade"; 'G' = 1; 'H' = 2; 'I' = 0; 'J' = 2; 'K' = "A"; 'L' = 8; 'M' = 5.25; 'N' = 1.285; 'O' = 15; 'P' = 7.5; 'Q' = 1.125; 'R' = 2.25; 'S' = 1.975; 'T' = 1.825; 'U' = 3.5; 'V' = 1.825; 'W' = 1.975; 'X' = -1; 'Y' = -1; 'Z' = 0.125; 'AA' = 0.9750000000000001; 'AB' = -1; 'AC' = -1; 'AD' = 0.9750000000000001 }
    },
    @{
        RowA = 279
        RowB = 280
        ValuesForRowA = @{ 'B' = 8106964; 'E' = "Radnicki Nis"; 'F' = "IMT Novi Belgrade"; 'G' = 0; 'H' = 0; 'I' = 0; 'J' = 0; 'K' = "D"; 'L' = 2.1; 'M' = 3.1; 'N' = 3.2; 'O' = 2.15; 'P' = 3.3; 'Q' = 2.9; 'R' = -0.25; 'S' = 1.925; 'T' = 1.875; 'U' = 2.5; 'V' = 1.9; 'W' = 1.9; 'X' = -1; 'Y' = 2.3; 'Z' = -1; 'AA' = -0.5; 'AB' = 0.4375; 'AC' = -1; 'AD' = 0.8999999999999999 }
        ValuesForRowB = @{ 'B' = 8106767; 'E' = "Javor Ivanjica"; 'F' = "FK Vozdovac"; 'G' = 1; 'H' = 1; 'I' = 1; 'J' = 1; 'K' = "D"; 'L' = 2.25; 'M' = 3; 'N' = 3; 'O' = 2.15; 'P' = 3; 'Q' = 3.1; 'R' = -0.25; 'S' = 1.9; 'T' = 1.9; 'U' = 2.25; 'V' = 1.9; 'W' = 1.9; 'X' = -1; 'Y' = 2; 'Z' = -1; 'AA' = -0.5; 'AB' = 0.45; 'AC' = -0.5; 'AD' = 0.45 }
    },
    @{
        RowA = 284
        RowB = 285
        ValuesForRowA = @{ 'B' = 8105026; 'E' = "FK Cukaricki"; 'F' = "Partizan Belgrade"; 'G' = 0; 'H' = 1; 'I' = 0; 'J' = 0; 'K' = "A"; 'L' = 2.75; 'M' = 3.2; 'N' = 2.3; 'O' = 2.45; 'P' = 3.6; 'Q' = 2.45; 'R' = 0; 'S' = 1.875; 'T' = 1.925; 'U' = 3; 'V' = 1.8; 'W' = 2; 'X' = -1; 'Y' = -1; 'Z' = 1.45; 'AA' = -1; 'AB' = 0.925; 'AC' = -1; 'AD' = 1 }
        ValuesForRowB = @{ 'B' = 8105865; 'E' = "Mladost Lucani"; 'F' = "FK Backa Topola"; 'G' = 2; 'H' = 0; 'I' = 1; 'J' = 0; 'K' = "H"; 'L' = 3.75; 'M' = 3.75; 'N' = 1.727; 'O' = 5.25; 'P' = 4.5; 'Q' = 1.42; 'R' = 1.25; 'S' = 1.9; 'T' = 1.9; 'U' = 3; 'V' = 1.9; 'W' = 1.9; 'X' = 4.25; 'Y' = -1; 'Z' = -1; 'AA' = 0.8999999999999999; 'AB' = -1; 'AC' = -1; 'AD' = 0.8999999999999999 }
    },
    @{
        RowA = 291
        RowB = 292
        ValuesForRowA = @{ 'B' = 8245725; 'E' = "Javor Ivanjica"; 'F' = "FK Zeleznicar Pancevo"; 'G' = 1; 'H' = 1; 'I' = 0; 'J' = 0; 'K' = "D"; 'L' = 2.375; 'M' = 2.875; 'N' = 2.9; 'O' = 3.8; 'P' = 2.15; 'Q' = 2.625; 'R' = 0; 'S' = 2.1; 'T' = 1.7; 'U' = 2.25; 'V' = 2.05; 'W' = 1.75; 'X' = -1; 'Y' = 1.15; 'Z' = -1; 'AA' = 0; 'AB' = 0; 'AC' = -0.5; 'AD' = 0.375 }
        ValuesForRowB = @{ 'B' = 8245726; 'E' = "Spartak Subotica"; 'F' = "IMT Novi Belgrade"; 'G' = 1; 'H' = 1; 'I' = 0; 'J' = 0; 'K' = "D"; 'L' = 4.4; 'M' = 2.75; 'N' = 1.909; 'O' = 4; 'P' = 2.875; 'Q' = 1.95; 'R' = 0.5; 'S' = 1.8; 'T' = 2; 'U' = 2.25; 'V' = 1.95; 'W' = 1.85; 'X' = -1; 'Y' = 1.875; 'Z' = -1; 'AA' = 0.8; 'AB' = -1; 'AC' = -0.5; 'AD' = 0.425 }
    }
)

foreach ($swap in $swaps) {
    foreach ($col in $swap.ValuesForRowA.Keys) {
        $ws.Range("$col$($swap.RowA)").Value = $swap.ValuesForRowA[$col]
    }
    foreach ($col in $swap.ValuesForRowB.Keys) {
        $ws.Range("$col$($swap.RowB)").Value = $swap.ValuesForRowB[$col]
    }
}
